# Update Deus Phase 1 / Phase 2 enemy stat rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed the new unique names first so they land in the shared-string table
# in the same order as the source workbook (Deus_Phase_1, Deus_Phase_2,
# Attack_01, Shot_01, Shot_02, Get_Hit, Dead_02).
$ws.Range("A13").Value = "Deus_Phase_1"
$ws.Range("A14").Value = "Deus_Phase_2"
$ws.Range("L13").Value = "Attack_01"
$ws.Range("M13").Value = "Shot_01"
$ws.Range("N13").Value = "Shot_02"
$ws.Range("O13").Value = "Get_Hit"
$ws.Range("P13").Value = "Dead_02"

# Row 13: Deus_Phase_1
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 30
$ws.Range("F13").Value = 500
$ws.Range("G13").Value = 250
$ws.Range("H13").Value = 300
$ws.Range("I13").Value = 7
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 3000

# Row 14: Deus_Phase_2
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 30
$ws.Range("F14").Value = 500
$ws.Range("G14").Value = 250
$ws.Range("H14").Value = 300
$ws.Range("I14").Value = 7
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = "Attack_01"
$ws.Range("M14").Value = "Shot_01"
$ws.Range("N14").Value = "Shot_02"
$ws.Range("O14").Value = "Get_Hit"
$ws.Range("P14").Value = "Dead_02"

# Sheet view adjustment: update the active selection like the author did
$ws.Range("G18").Select()
